$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("APR-2021")

# ---------------------------------------------------------------------------
# Row 17 (existing row): update the "Application" (C) and "Task" (D) columns
# and grow the row height to fit the extra wrapped line of text.
# ---------------------------------------------------------------------------
$ws.Cells.Item(17, 3).Value = "Muji store , B2B & nMVAR "
$ws.Cells.Item(17, 4).Value = "Sanity testing on B2C app, QMVAR site, GSS site and Hayaai site. `nRegression testing and Retesting on Muji store application`nRegression testing and Retesting on nMVAR application`nRegression testing and Retesting on B2B application"
$ws.Rows.Item(17).RowHeight = 60

# ---------------------------------------------------------------------------
# Row 18: "Holiday" (Saturday 17-Apr-2021) - same layout as row 14/15
# ---------------------------------------------------------------------------
$ws.Range("A14:G14").Copy($ws.Range("A18:G18"))
$ws.Cells.Item(18, 1).Value = 17
$ws.Cells.Item(18, 2).Value = 44303

# ---------------------------------------------------------------------------
# Row 19: "Holiday" (Sunday 18-Apr-2021)
# ---------------------------------------------------------------------------
$ws.Range("A14:G14").Copy($ws.Range("A19:G19"))
$ws.Cells.Item(19, 1).Value = 18
$ws.Cells.Item(19, 2).Value = 44304

# ---------------------------------------------------------------------------
# Row 20: Monday 19-Apr-2021 - filled task row, same layout as row 9/10
# ---------------------------------------------------------------------------
$ws.Range("A9:G9").Copy($ws.Range("A20:G20"))
$ws.Cells.Item(20, 1).Value = 19
$ws.Cells.Item(20, 2).Value = 44305
$ws.Cells.Item(20, 3).Value = "Muji store & B2B"
$ws.Cells.Item(20, 4).Value = "Sanity testing on B2C app, QMVAR site, GSS site and Hayaai site. `nRegression testing and Retesting on Muji store application`nRegression testing and Retesting on B2B application"
$ws.Rows.Item(20).RowHeight = 45

# ---------------------------------------------------------------------------
# Rows 21-31: future / unfilled days (20-Apr-2021 .. 30-Apr-2021)
# Layout: No./Date filled in, Application/Task/%/Status/Comments left blank
# but still carrying the same cell formatting as the rest of the table.
# ---------------------------------------------------------------------------
$ws.Range("A11:C11").Copy($ws.Range("A21:C21"))
$ws.Range("E11:G11").Copy($ws.Range("E21:G21"))
$ws.Range("C8").Copy($ws.Range("D21"))
$ws.Range("D21").ClearContents()

$blankDates = 44306, 44307, 44308, 44309, 44310, 44311, 44312, 44313, 44314, 44315, 44316
$blankNo = 20, 21, 22, 23, 24, 25, 26, 27, 28, 29, 30

for ($i = 0; $i -lt $blankDates.Length; $i++) {
    $r = 21 + $i
    if ($r -gt 21) {
        $ws.Range("A21:G21").Copy($ws.Range("A" + $r + ":G" + $r))
    }
    $ws.Cells.Item($r, 1).Value = $blankNo[$i]
    $ws.Cells.Item($r, 2).Value = $blankDates[$i]
}
